$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.477.05"
$ws.Range("E2").Value = "  +2.36%  "

# Row 3
$ws.Range("D3").Value = "2.510.91"
$ws.Range("E3").Value = "  +0.66%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "597.11"
$ws.Range("E5").Value = "  +1.65%  "

# Row 6
$ws.Range("D6").Value = "176.10"
$ws.Range("E6").Value = "  -0.70%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  +0.55%  "

# Row 9
$ws.Range("D9").Value = "2.509.64"
$ws.Range("E9").Value = "  +0.64%  "

# Row 10
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  +10.91%  "

# Row 11
$ws.Range("E11").Value = "  -0.48%  "

# Row 12
$ws.Range("E12").Value = "  +0.62%  "

# Row 13
$ws.Range("D13").Value = "5.00"
$ws.Range("E13").Value = "  +1.37%  "

# Row 14
$ws.Range("D14").Value = "2.968.87"
$ws.Range("E14").Value = "  +0.66%  "

# Row 15
$ws.Range("D15").Value = "25.93"
$ws.Range("E15").Value = "  +0.87%  "

# Row 16
$ws.Range("D16").Value = "69.424.70"
$ws.Range("E16").Value = "  +2.51%  "

# Row 17
$ws.Range("E17").Value = "  +3.03%  "

# Row 18
$ws.Range("D18").Value = "2.514.79"
$ws.Range("E18").Value = "  +0.20%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  +1.17%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "362.21"
$ws.Range("E20").Value = "  +3.15%  "

# Row 21
$ws.Range("D21").Value = "11.02"
$ws.Range("E21").Value = "  +0.30%  "

# Row 22
$ws.Range("D22").Value = "4.05"
$ws.Range("E22").Value = "  -1.68%  "

# Row 23
$ws.Range("E23").Value = "  -0.14%  "

# Row 24
$ws.Range("D24").Value = "70.50"
$ws.Range("E24").Value = "  -0.53%  "

# Row 25
$ws.Range("D25").Value = "4.22"
$ws.Range("E25").Value = "  -2.11%  "

# Row 26
$ws.Range("E26").Value = "  -0.99%  "

# Row 27
$ws.Range("E27").Value = "  -3.14%  "

# Row 28
$ws.Range("D28").Value = "2.627.73"
$ws.Range("E28").Value = "  +0.30%  "

# Row 29
$ws.Range("E29").Value = "  +0.58%  "

# Row 30
$ws.Range("D30").Value = "512.44"
$ws.Range("E30").Value = "  +1.35%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0893"
$ws.Range("E31").Value = "  -1.70%  "

# Row 32
$ws.Range("D32").Value = "7.76"
$ws.Range("E32").Value = "  -1.14%  "

# Row 33
$ws.Range("D33").Value = "1.24"
$ws.Range("E33").Value = "  -1.81%  "

# Row 34
$ws.Range("E34").Value = "  +0.57%  "

# Row 35
$ws.Range("E35").Value = "  +0.03%  "

# Row 36
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "162.80"
$ws.Range("E36").Value = "  -0.10%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.119"
$ws.Range("E37").Value = "  -1.95%  "

# Row 38
$ws.Range("D38").Value = "18.73"
$ws.Range("E38").Value = "  +2.15%  "

# Row 39
$ws.Range("E39").Value = "  +1.23%  "

# Row 40
$ws.Range("E40").Value = "  -1.23%  "

# Row 41
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("E42").Value = "  -1.56%  "

# Row 43
$ws.Range("D43").Value = "4.79"
$ws.Range("E43").Value = "  -1.55%  "

# Row 44
$ws.Range("D44").Value = "0.320"
$ws.Range("E44").Value = "  -2.45%  "

# Row 45
$ws.Range("E45").Value = "  -3.88%  "

# Row 46
$ws.Range("D46").Value = "38.78"
$ws.Range("E46").Value = "  -0.58%  "

# Row 47
$ws.Range("D47").Value = "149.87"
$ws.Range("E47").Value = "  +3.45%  "

# Row 48
$ws.Range("E48").Value = "  +1.59%  "

# Row 49
$ws.Range("E49").Value = "  +0.10%  "

# Row 50
$ws.Range("D50").Value = "0.0₆0252"
$ws.Range("E50").Value = "  -0.94%  "

# Row 51
$ws.Range("E51").Value = "  -0.80%  "
